# Update "想去人数" (want-to-go count) values in column F for the two
# sheets that carry the full event data: "展览" and "全部类型".
# Rows (by sheet row number) and their new F-column values:
#   F2: 8347 -> 8350
#   F3: 7799 -> 7804
#   F5: 190  -> 191
#   F10: 165 -> 166
#   F14: 1343 -> 1344
#   F19: 123 -> 124

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8350
    3  = 7804
    5  = 191
    10 = 166
    14 = 1344
    19 = 124
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
